$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.302.68"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.588.54"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'209.92"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "'19.46"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.812.02"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.635.31"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.07"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'0.518"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "'64.31"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "26.312.66"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "'7.47"
$ws.Range("E19").Value = "  +5.46%  "
$ws.Range("D20").Value = "'210.97"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'144.83"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "'15.23"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "1.301.70"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.44"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.611"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").Value = "  -11.15%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'5.61"
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").Value = "'62.32"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").Value = "1.724.37"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").Value = "'87.76"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("E48").Value = "  -5.36%  "
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").Value = "'0.0981"
$ws.Range("E50").Value = "  -4.62%  "
$ws.Range("E51").Value = "  -0.29%  "
